# Update economic factor impact by adjusting weightings for different sectors
# This adjusts specific cells on the "Impact Values (1-3)" sheet, and the
# corresponding derived cells on the "Weights (Impact x Importance)" sheet
# (since those rows have an Importance multiplier of 1, the raw impact value
# and the weight value are numerically identical for these cells).

$wb = $excel.ActiveWorkbook

$impactSheet = $wb.Worksheets.Item("Impact Values (1-3)")
$weightsSheet = $wb.Worksheets.Item("Weights (Impact x Importance)")

# Cell -> new value pairs shared by both sheets
$updates = @{
    "B6"  = 3
    "B7"  = 3
    "B10" = 3
    "E10" = 3
    "G10" = 3
    "H10" = 3
    "E11" = 2
    "F11" = 2
    "G11" = 3
    "H11" = 2
    "B13" = 2
    "C13" = 2
    "F13" = 2
    "J13" = 2
}

foreach ($sheet in @($impactSheet, $weightsSheet)) {
    foreach ($addr in $updates.Keys) {
        $sheet.Range($addr).Value = $updates[$addr]
    }
}
